# "add scrollable feature to the visualization applet"
# - widen column C and zoom the sheet view in (170%) so the long notes are readable
# - move selection to A13
# - strip the leading checkmark glyphs ("✔；" / "✔; ") from the existing D-column notes
# - clear the now-empty D9 note
# - append two new log rows (817Train_ctm200ep_rd / 818Train_ctm200ep_rd)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- widen column C (65.33 -> ~83.33 chars) ---
$ws.Columns.Item(3).ColumnWidth = 82.46

# --- zoom in on the sheet (applet becomes scrollable at higher zoom) and move selection ---
$win = $excel.ActiveWindow
$win.Zoom = 170
$ws.Range("A13").Select()

# --- clean up the D-column notes: drop the leading checkmark markers ---
$ws.Range("D3").Value = 'have some abnormity: Rd_Smooth image sets give large error'
$ws.Range("D4").Value = 'looks good'
$ws.Range("D5").Value = 'trained on the wrong image set ✖️'
$ws.Range("D6").Value = 'Rd_Smooth image still give poor results than others'
$ws.Range("D7").Value = '
Weird patterns, reach a flat plateau in messiness prediction; 
might be the problem of loss function???
Prediction Visualization shows better recognition of complexity within a central area; '

# D9 note removed entirely
$ws.Range("D9").ClearContents()

$ws.Range("D10").Value = 'R: recognizes overall roundness better; 
still bad at recognize the color-smooth round shapes;
M: personally wouldn’t agree with the prediction 
'
$ws.Range("D11").Value = ' recognizes the inner noise better! 
'

# row 11 no longer needs the extra-tall wrap height
$ws.Rows.Item(11).RowHeight = 32

# --- append new rows for the two newest training runs ---
$ws.Range("A15").Value = "817Train_ctm200ep_rd"
$ws.Range("B15").Value = "As above "
$ws.Range("C15").Value = 'Refactored processing code so that the stroke color and background color differ by at least 70; 
removed images with blur=3/4; '
$ws.Range("C15").HorizontalAlignment = -4108
$ws.Range("C15").VerticalAlignment = -4108
$ws.Range("C15").WrapText = $true
$ws.Range("D15").Value = 'lower test score -> less data images?;
Recognizes the inperceptiable gradients less, especially with "Images/evo_art_test/Image-1.jpg";
Gives much lower score to blurred circles, shown by "evo_art_test/image-49.jpg"
'
$ws.Range("D15").HorizontalAlignment = -4131
$ws.Range("D15").VerticalAlignment = -4108
$ws.Range("D15").WrapText = $true
$ws.Rows.Item(15).RowHeight = 154

$ws.Range("A16").Value = "818Train_ctm200ep_rd"
$ws.Range("B16").Value = "As above "
$ws.Range("C16").Value = "added images with blur=3/4 but with heavier stroke weights;"
